# Auto-generated edit script applying the Kujata_Profits (Leve profit) data refresh
# across all 8 job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H28").Value = 3142.6924
$ws.Range("I28").Value = 2623.182
$ws.Range("J28").Value = 6000
$ws.Range("K28").Value = 2623.182
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = -2138.182
$ws.Range("N28").Value = -6970

$ws.Range("H29").Value = 900
$ws.Range("J29").Value = 2600
$ws.Range("L29").Value = 7800
$ws.Range("N29").Value = -8362

$ws.Range("H38").Value = 187.4
$ws.Range("I38").Value = 187.4
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 562.2
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -190.2
$ws.Range("N38").ClearContents()

$ws.Range("H40").Value = 1350.5

$ws.Range("H43").Value = 18058154
$ws.Range("I43").Value = 1022.75
$ws.Range("J43").Value = 30096244
$ws.Range("K43").Value = 1022.75
$ws.Range("L43").Value = 30096244
$ws.Range("M43").Value = -953.75
$ws.Range("N43").Value = -30096382

$ws.Range("H58").Value = 2532.5881
$ws.Range("J58").Value = 3456.7273
$ws.Range("L58").Value = 10370.1819
$ws.Range("N58").Value = -10670.1819

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H135").Value = 1359.697
$ws.Range("I135").Value = 606.0714
$ws.Range("J135").Value = 5580
$ws.Range("K135").Value = 5454.6426
$ws.Range("L135").Value = 50220
$ws.Range("M135").Value = -2919.6426
$ws.Range("N135").Value = -55290

$ws.Range("H137").Value = 3138.1667
$ws.Range("I137").Value = 3037
$ws.Range("J137").Value = 3401.2
$ws.Range("K137").Value = 9111
$ws.Range("L137").Value = 10203.6
$ws.Range("M137").Value = -6561
$ws.Range("N137").Value = -15303.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12448.5
$ws.Range("I32").Value = 9558.120000000001
$ws.Range("J32").Value = 17265.8
$ws.Range("K32").Value = 9558.120000000001
$ws.Range("L32").Value = 17265.8
$ws.Range("M32").Value = -9271.120000000001
$ws.Range("N32").Value = -17839.8

$ws.Range("H61").Value = 83334980
$ws.Range("I61").Value = 125000904
$ws.Range("K61").Value = 125000904
$ws.Range("M61").Value = -125000692

$ws.Range("H74").Value = 2254.318
$ws.Range("I74").Value = 1147.2307
$ws.Range("J74").Value = 3853.4443
$ws.Range("K74").Value = 1147.2307
$ws.Range("L74").Value = 3853.4443
$ws.Range("M74").Value = -273.2307000000001
$ws.Range("N74").Value = -5601.4443

$ws.Range("H77").Value = 2254.318
$ws.Range("I77").Value = 1147.2307
$ws.Range("J77").Value = 3853.4443
$ws.Range("K77").Value = 5736.1535
$ws.Range("L77").Value = 19267.2215
$ws.Range("M77").Value = -1368.1535
$ws.Range("N77").Value = -28003.2215

$ws.Range("H122").Value = 3531.4
$ws.Range("I122").Value = 2912.2856
$ws.Range("J122").Value = 4976
$ws.Range("K122").Value = 8736.856800000001
$ws.Range("L122").Value = 14928
$ws.Range("M122").Value = -6286.856800000001
$ws.Range("N122").Value = -19828

$ws.Range("H136").Value = 83334980
$ws.Range("I136").Value = 125000904
$ws.Range("K136").Value = 375002712
$ws.Range("M136").Value = -375000162

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 43482096
$ws.Range("I86").Value = 47622816
$ws.Range("K86").Value = 47622816
$ws.Range("M86").Value = -47621693

$ws.Range("H89").Value = 43482096
$ws.Range("I89").Value = 47622816
$ws.Range("K89").Value = 238114080
$ws.Range("M89").Value = -238108464

$ws.Range("H134").Value = 10460.272
$ws.Range("I134").Value = 1135.375
$ws.Range("J134").Value = 35326.668
$ws.Range("K134").Value = 3406.125
$ws.Range("L134").Value = 105980.004
$ws.Range("M134").Value = -871.125
$ws.Range("N134").Value = -111050.004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6261
$ws.Range("I58").Value = 916.75
$ws.Range("J58").Value = 9549.77
$ws.Range("K58").Value = 916.75
$ws.Range("L58").Value = 9549.77
$ws.Range("M58").Value = -713.75
$ws.Range("N58").Value = -9955.77

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
$ws.Range("M122").Value = -2050

$ws.Range("H134").Value = 10418046
$ws.Range("I134").Value = 1429.6129
$ws.Range("J134").Value = 29413052
$ws.Range("K134").Value = 4288.8387
$ws.Range("L134").Value = 88239156
$ws.Range("M134").Value = -1753.8387
$ws.Range("N134").Value = -88244226

$ws.Range("H136").Value = 6261
$ws.Range("I136").Value = 916.75
$ws.Range("J136").Value = 9549.77
$ws.Range("K136").Value = 2750.25
$ws.Range("L136").Value = 28649.31
$ws.Range("M136").Value = -200.25
$ws.Range("N136").Value = -33749.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1605.7273
$ws.Range("I5").Value = 1470.0526
$ws.Range("J5").Value = 2465
$ws.Range("K5").Value = 4410.1578
$ws.Range("L5").Value = 7395
$ws.Range("M5").Value = -4298.1578
$ws.Range("N5").Value = -7619

$ws.Range("H21").Value = 999.5
$ws.Range("J21").Value = 1799
$ws.Range("L21").Value = 5397
$ws.Range("N21").Value = -5743

$ws.Range("H68").Value = 936.5833
$ws.Range("I68").Value = 1123.75
$ws.Range("J68").Value = 843
$ws.Range("K68").Value = 3371.25
$ws.Range("L68").Value = 2529
$ws.Range("M68").Value = -2560.25
$ws.Range("N68").Value = -4151

$ws.Range("H71").Value = 936.5833
$ws.Range("I71").Value = 1123.75
$ws.Range("J71").Value = 843
$ws.Range("K71").Value = 10113.75
$ws.Range("L71").Value = 7587
$ws.Range("M71").Value = -6057.75
$ws.Range("N71").Value = -15699

$ws.Range("H135").Value = 1605.7273
$ws.Range("I135").Value = 1470.0526
$ws.Range("J135").Value = 2465
$ws.Range("K135").Value = 13230.4734
$ws.Range("L135").Value = 22185
$ws.Range("M135").Value = -10695.4734
$ws.Range("N135").Value = -27255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1259.1
$ws.Range("I113").Value = 1351.8334
$ws.Range("J113").Value = 1120
$ws.Range("K113").Value = 1351.8334
$ws.Range("L113").Value = 1120
$ws.Range("M113").Value = 818.1666
$ws.Range("N113").Value = -5460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1184
$ws.Range("J7").Value = 1947.5
$ws.Range("L7").Value = 1947.5
$ws.Range("N7").Value = -2171.5

$ws.Range("H40").Value = 2175.6843
$ws.Range("J40").Value = 2575.25
$ws.Range("L40").Value = 2575.25
$ws.Range("N40").Value = -2847.25

$ws.Range("H126").Value = 1184
$ws.Range("J126").Value = 1947.5
$ws.Range("L126").Value = 5842.5
$ws.Range("N126").Value = -10782.5

$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws.Range("H136").Value = 2030
$ws.Range("I136").Value = 1922.2222
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5766.6666
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3216.6666
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1347.566
$ws.Range("I136").Value = 511.24243
$ws.Range("K136").Value = 1533.72729
$ws.Range("M136").Value = 1016.27271
